# Add "Parent Activities" tab entries to the error/string list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data (row, A id, B category, C message). Note: the shared-string table
# records "Total" before "contribution(s)" (the original author must have
# entered/edited row 264 ahead of row 263), so we write the cell values in
# that same order below even though the rows are populated out of visual
# sequence; the final cell placement still ends up correct.
$ws.Range("A260").Value = 1005
$ws.Range("B260").Value = "Parent"
$ws.Range("C260").Value = "Password doesn't match."
$ws.Range("D260").Value = "Message"

$ws.Range("A261").Value = 1006
$ws.Range("B261").Value = "Parent"
$ws.Range("C261").Value = "Password has to be longer than 5 characters."
$ws.Range("D261").Value = "Message"

$ws.Range("A262").Value = 1007
$ws.Range("B262").Value = "Parent"
$ws.Range("C262").Value = "PARENT ACTIVITIES"
$ws.Range("D262").Value = "Message"

$ws.Range("A264").Value = 1009
$ws.Range("B264").Value = "Parent"
$ws.Range("C264").Value = "Total"
$ws.Range("D264").Value = "Message"

$ws.Range("A263").Value = 1008
$ws.Range("B263").Value = "Parent"
$ws.Range("C263").Value = "contribution(s)"
$ws.Range("D263").Value = "Message"

$ws.Range("A265").Value = 1010
$ws.Range("B265").Value = "Parent"
$ws.Range("C265").Value = "since"
$ws.Range("D265").Value = "Message"

$ws.Range("A266").Value = 1011
$ws.Range("B266").Value = "Parent"
$ws.Range("C266").Value = "Recent 10 Posts"
$ws.Range("D266").Value = "Message"

$lastRow = 266

$ws.Range("C$lastRow").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 251
$win.ScrollColumn = 1
